$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header AH1 (new shared string "20-jul")
$ws.Range("AH1").Value = "20-jul"

# Fill AH2:AH18 with new data column values
$ws.Range("AH2").Value = 0
$ws.Range("AH3").Value = 8.3892039993575054
$ws.Range("AH4").Value = 12.958977668816591
$ws.Range("AH5").Value = 24.985887282395023
$ws.Range("AH6").Value = 0
$ws.Range("AH7").Value = 9.0259672675828213
$ws.Range("AH8").Value = 10.203274628609107
$ws.Range("AH9").Value = 18.994987121083639
$ws.Range("AH10").Value = 22.04276782764936
$ws.Range("AH11").Value = 13.215969599109227
$ws.Range("AH12").Value = 0
$ws.Range("AH13").Value = 9.0379072696581595
$ws.Range("AH14").Value = 0
$ws.Range("AH15").Value = 0
$ws.Range("AH16").Value = 12.5826256717438
$ws.Range("AH17").Value = 0
$ws.Range("AH18").Value = 0

# Update active selection to match the authored state
[void]$ws.Range("AH8").Select()

